$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold font + border, style index 1) from H1 to the
# two new header cells before setting their text, mirroring existing headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$values = @{
    2 = 5
    3 = 8
    4 = 8
    5 = 7
    6 = 6
    7 = 7
    8 = 6
    9 = 5
    10 = 5
    11 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("I$row").Value = $values[$row]
    $ws.Range("J$row").Value = $values[$row]
}
